$d = $word.ActiveDocument

# Locate the paragraph that follows the insertion point:
# "After this line there should be stars." (a FirstParagraph-styled
# paragraph that immediately follows the BlockText paragraph ending in
# "...and formatting.")
$target = $null
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "After this line there should be stars.*") {
        $target = $p
        $targetIndex = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the 'After this line there should be stars.' paragraph."
}

# The insertion point is the end of the previous paragraph (the one
# ending in "...and formatting."), i.e. just before its own paragraph
# mark, so that InsertXML adds whole new paragraphs in between without
# swallowing either neighbor.
$prev = $d.Paragraphs.Item($targetIndex - 1)
$insertPos = $prev.Range.End - 1

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xmlDouble = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val='BlockText'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>“These literal double curly quotes, used where smart</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>quotes gets it wrong, curl the right way even though</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>they’re on different lines.”</w:t></w:r>" +
    "</w:p>"

$xmlSingle = "<w:p $wNs>" +
    "<w:pPr><w:pStyle w:val='BlockText'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>‘These literal single curly quotes, used where smart</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>quotes gets it wrong, curl the right way even though</w:t></w:r>" +
    "<w:r><w:br/></w:r>" +
    "<w:r><w:t xml:space='preserve'>they’re on different lines.’</w:t></w:r>" +
    "</w:p>"

# Insert the "double quotes" paragraph first.
$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xmlDouble)

# Re-find the "After this line..." paragraph, and insert the "single
# quotes" paragraph directly before it (i.e. right after the paragraph
# we just inserted).
$target = $null
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "After this line there should be stars.*") {
        $target = $p
        $targetIndex = $i
        break
    }
}
$prev = $d.Paragraphs.Item($targetIndex - 1)
$insertPos = $prev.Range.End - 1

$r = $d.Range($insertPos, $insertPos)
$r.InsertXML($xmlSingle)

Write-Host "done"
